$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 101, pushing the existing rows 101-134
# down to 102-135 (this also extends the used range to A1:R135, matching
# the inherited date-style on column D from the row above).
$ws.Rows.Item(101).Insert()

# Populate the newly inserted row 101 with the new weekly record.
$ws.Cells.Item(101, 1).Value = 8
$ws.Cells.Item(101, 2).Value = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(101, 3).Value = 'Coquimbo'
$ws.Cells.Item(101, 4).Value = '2021-11-24'
$ws.Cells.Item(101, 5).Value = 4
$ws.Cells.Item(101, 6).Value = 100112037
$ws.Cells.Item(101, 7).Value = 'Cebollín'
$ws.Cells.Item(101, 8).Value = 'Sin especificar'
$ws.Cells.Item(101, 9).Value = 'Primera'
$ws.Cells.Item(101, 10).Value = 3200
$ws.Cells.Item(101, 11).Value = 900
$ws.Cells.Item(101, 12).Value = 1000
$ws.Cells.Item(101, 13).Value = 950
$ws.Cells.Item(101, 14).Value = '$/paquete 6 unidades'
$ws.Cells.Item(101, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(101, 16).Value = 158
$ws.Cells.Item(101, 17).Value = 6
$ws.Cells.Item(101, 18).Value = 'Hortaliza'
